$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.263.58"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.690.35"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5263"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.007"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2697"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06452"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07477"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.567"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.682.45"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5862"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008545"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "26.306.06"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.976"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.241"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.688"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1238"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06642"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.357"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.331"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.594"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.563"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.672"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.030"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6223"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.391"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.706"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.301"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01623"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("D40").Value = "1.103.07"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8878"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.016"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("D44").Value = "1.837.66"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000110"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.166"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05263"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4299"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.060"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.88%  "
